$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text (matching source inlineStr/text cells),
# so numeric-looking strings like "19.27" are not auto-converted to numbers.
$targetCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "E6", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "E23", "E24", "D25", "E25", "E26", "E27", "D28", "E28", "D29", "E29", "E30", "E31", "E32", "E33", "D34", "E34", "E35", "B36", "C36", "D36", "E36", "B37", "C37", "D37", "E37", "D38", "E38", "E39", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "E50", "E51")
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "26.659.62"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.643.60"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "215.37"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "0.0628"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").Value = "19.27"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.873.67"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.19"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.633.57"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "0.530"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "65.27"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "26.686.70"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "216.68"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  +15.04%  "
$ws.Range("D25").Value = "146.24"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("D29").Value = "15.75"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").Value = "1.271.96"
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0180"
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "0.817"
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "1.783.51"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").Value = "91.63"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "59.54"
$ws.Range("E46").Value = "  +8.17%  "
$ws.Range("D47").Value = "1.60"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "7.82"
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("E51").Value = "  -0.52%  "
